$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Marking" row (row 11): Right marks 4 -> 5, Wrong marks -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Update "Total" row (row 12): Right total 80 -> 100, Wrong total -8 -> -9.6
$ws.Range("B12").Value = 100
$ws.Range("C12").Value = -9.6

# Update the score string in E12: "72/112" -> "90.4/140"
$ws.Range("E12").Value = "90.4/140"
